$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (scenario2) mirroring scenario1's structure
$data = @(
    @{ A="scenario2"; B="LIVEHTA Automation - Test_NonOncology_Automation_3"; C="LIVEHTA Automation - Test_NonOncology_Automation_3_radio_button"; D="Clinical"; E="Clinical_radio_button"; F="sub_pop_section1"; G="sub_pop_section1_checkbox"; H="sub_pop_section"; I="StandardExcelReport-LIVEHTA Automation - Test_NonOncology_Automation_3-Clinical-2023_"; J=7 },
    @{ A="scenario2"; F="intervention_section4"; G="intervention_section4_checkbox"; H="intervention_section"; I="CompleteExcelReport-LIVEHTA Automation - Test_NonOncology_Automation_3-Clinical-2023_"; J=3 },
    @{ A="scenario2"; F="study_design_section1"; G="study_design_section1_checkbox"; H="study_design_section"; J=0 },
    @{ A="scenario2"; J=1 },
    @{ A="scenario2"; J=1 },
    @{ A="scenario2"; J=0 },
    @{ A="scenario2"; J=4 }
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $data[$i]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}

$ws.Range("A16:XFD16").Select()
